$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in new row 11: email = "lalalala", condition = "emailSalah"
$ws.Range("A11").Value = "lalalala"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "emailSalah"

# Copy style from row 10 (A10/F10) onto row 11 so formatting matches the other data rows
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("B10:E10").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122) | Out-Null

# Add hyperlink on A11 like the other email cells
$ws.Hyperlinks.Add($ws.Range("A11"), "mailto:rizkariz20@gmail.com", "", "rizkariz20@gmail.com", "lalalala")

# Update the active selection to A12 (reflects where the user clicked last)
$ws.Range("A12").Select()
